# Apply the commit: "fixed LBNRIND, docs and delete check sites"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo LBNDIND -> LBNRIND in the F1 header cell
$ws.Range("F1").Value = "LBNRIND"

# Update the selection to K5 (single cell, no longer A1:H1 range)
$ws.Range("K5").Select()
